# Knowledge base references polishment.
# - Re-labels the reference links on Sheet2 with section headers.
# - Re-points the old pss.uvm.edu link at a Cactus nursery pH source and
#   adds a trailing space to the almanac link.
# - Adds a new "Cactus" row to the KB lookup table (sheet1, row 113).
# - Activates Sheet2 as the visible tab with a fresh selection.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet2 content edits (order matters for shared-string allocation) ---

# 1) Move the almanac link text (was B5) down to B8, appending a trailing
#    space to its display text, then clear the old B5 cell.
$ws2.Range("B8").Value = "https://www.almanac.com/plant-ph# "
$ws2.Range("B5").ClearContents()

# 2) Move the old pss.uvm.edu link (was B7) down to B12 as the new Cactus
#    nursery pH link, then clear out the stale B7 cell before re-using it.
$ws2.Range("B12").Value = "https://www.cactusnursery.co.uk/ph.htm "
$ws2.Range("B7").ClearContents()

# 3) New section-header labels.
$ws2.Range("B3").Value = "For row 2 to 18"
$ws1.Range("A113").Value = "Cactus"
$ws2.Range("B11").Value = "For Cactus (Row 113)"
$ws2.Range("B7").Value = "For row 19 to 112"

# --- Finish the new KB row (sheet1 row 113) ---
$ws1.Range("B113").Value = 5
$ws1.Range("C113").Value = 6.5

# --- Hyperlinks + Hyperlink style on Sheet2 ---
$ws2.Hyperlinks.Add($ws2.Range("B4"), "https://www.cropnutrition.com/efu-soil-ph")
$ws2.Hyperlinks.Add($ws2.Range("B8"), "https://www.almanac.com/plant-ph#", " ")
$ws2.Hyperlinks.Add($ws2.Range("B12"), "https://www.cactusnursery.co.uk/ph.htm")

# B10 carries the Hyperlink style but no value/hyperlink object of its own.
$ws2.Range("B10").Style = "Hyperlink"

# --- View / selection state ---
$ws1.Select()
$ws1.Range("C125").Select()
$ws2.Select()
$ws2.Range("G10").Select()
